# Add a new "2022-Q1" sheet (holdings detail) between "2021-Q4" and "总计",
# and record its totals in the "总计" summary sheet.

$wb = $excel.ActiveWorkbook
$q4Sheet = $wb.Worksheets.Item("2021-Q4")

# --- Create the new "2022-Q1" sheet right after "2021-Q4" ---
$newSheet = $wb.Worksheets.Add($null, $q4Sheet)
$newSheet.Name = "2022-Q1"

# Re-use the existing bordered/bold style from the "2021-Q4" sheet instead of
# building new style entries from scratch.
$q4Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$q4Sheet.Range("A2").Copy()
$newSheet.Range("A2:A5").PasteSpecial(-4122)

# Headers
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data rows
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'920003"
$newSheet.Range("C2").Value = "中金新锐股票A"
$newSheet.Range("D2").Value = "'24.64"
$newSheet.Range("E2").Value = "'92.76"
$newSheet.Range("F2").Value = "'2.94"
$newSheet.Range("G2").Value = "'0.7244"
$newSheet.Range("H2").Value = 7

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'920923"
$newSheet.Range("C3").Value = "中金新锐股票C"
$newSheet.Range("D3").Value = "'3.94"
$newSheet.Range("E3").Value = "'92.76"
$newSheet.Range("F3").Value = "'2.94"
$newSheet.Range("G3").Value = "'0.1158"
$newSheet.Range("H3").Value = 7

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "'000926"
$newSheet.Range("C4").Value = "中信建投睿信灵活配置混合A"
$newSheet.Range("D4").Value = "'0.13"
$newSheet.Range("E4").Value = "'40.35"
$newSheet.Range("F4").Value = "'3.77"
$newSheet.Range("G4").Value = "'0.0049"
$newSheet.Range("H4").Value = 5

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "'004676"
$newSheet.Range("C5").Value = "中信建投睿信灵活配置混合C"
$newSheet.Range("D5").Value = "'0.02"
$newSheet.Range("E5").Value = "'40.35"
$newSheet.Range("F5").Value = "'3.77"
$newSheet.Range("G5").Value = "'0.0008"
$newSheet.Range("H5").Value = 5

# --- Update the "总计" sheet: insert a new row for 2022-Q1 before the
#     existing 2021-Q4 row, pushing it down. ---
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# Style the newly inserted A2 the same way as the (now shifted) A3.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 4
$totalSheet.Range("D2").Value = 0.85

$totalSheet.Range("A3").Value = 1
